$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.727.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "'2.338.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'238.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'0.664"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").Value = "'71.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.18%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.00%  "
$ws.Range("D10").Value = "'0.0989"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "'57.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'32.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.07%  "
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "'7.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.25%  "
$ws.Range("D15").Value = "'2.684.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "'16.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.45%  "
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("D18").Value = "'2.334.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'43.639.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'77.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'6.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'250.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  +6.08%  "
$ws.Range("D26").Value = "'3.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.71%  "
$ws.Range("D27").Value = "'2.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("E28").Value = "  -6.56%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'176.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'22.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").Value = "'0.0735"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("D35").Value = "'5.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.31%  "
$ws.Range("D36").Value = "'5.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").Value = "'5.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +30.71%  "
$ws.Range("D39").Value = "'6.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D41").Value = "'0.0270"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").Value = "'66.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.78%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'9.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.107"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").Value = "'18.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.46%  "
$ws.Range("D46").Value = "'0.194"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.65%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("E51").Value = "  -3.84%  "
